$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-CellText $ws.Range('D2') '94.556.70'
Set-CellText $ws.Range('E2') '  -2.99%  '
Set-CellText $ws.Range('D3') '3.442.03'
Set-CellText $ws.Range('E3') '  +2.99%  '
Set-CellText $ws.Range('E4') '  +0.01%  '
Set-CellText $ws.Range('D5') '238.39'
Set-CellText $ws.Range('E5') '  -4.88%  '
Set-CellText $ws.Range('D6') '644.22'
Set-CellText $ws.Range('E6') '  -1.82%  '
Set-CellText $ws.Range('D7') '1.44'
Set-CellText $ws.Range('E7') '  +4.07%  '
Set-CellText $ws.Range('D8') '0.406'
Set-CellText $ws.Range('E8') '  -3.42%  '
Set-CellText $ws.Range('E9') '  +0.09%  '
Set-CellText $ws.Range('D10') '0.974'
Set-CellText $ws.Range('E10') '  -1.83%  '
Set-CellText $ws.Range('D11') '3.441.33'
Set-CellText $ws.Range('E11') '  +3.08%  '
Set-CellText $ws.Range('D12') '42.67'
Set-CellText $ws.Range('E12') '  +4.50%  '
Set-CellText $ws.Range('D13') '0.199'
Set-CellText $ws.Range('E13') '  -4.42%  '
Set-CellText $ws.Range('D14') '6.23'
Set-CellText $ws.Range('E14') '  +2.41%  '
Set-CellText $ws.Range('D15') '94.251.10'
Set-CellText $ws.Range('E15') '  -3.09%  '
Set-CellText $ws.Range('D16') '4.082.35'
Set-CellText $ws.Range('E16') '  +2.93%  '
Set-CellText $ws.Range('D17') '0.0000252'
Set-CellText $ws.Range('E17') '  +0.03%  '
Set-CellText $ws.Range('D18') '8.40'
Set-CellText $ws.Range('E18') '  -1.62%  '
Set-CellText $ws.Range('D19') '3.452.31'
Set-CellText $ws.Range('E19') '  +3.79%  '
Set-CellText $ws.Range('D20') '17.70'
Set-CellText $ws.Range('E20') '  +0.88%  '
Set-CellText $ws.Range('D21') '11.43'
Set-CellText $ws.Range('E21') '  +6.37%  '
Set-CellText $ws.Range('D22') '0.503'
Set-CellText $ws.Range('E22') '  -3.50%  '
Set-CellText $ws.Range('D23') '502.40'
Set-CellText $ws.Range('E23') '  -0.60%  '
Set-CellText $ws.Range('E24') '  -2.86%  '
Set-CellText $ws.Range('D25') '0.0000194'
Set-CellText $ws.Range('E25') '  -2.61%  '
Set-CellText $ws.Range('D26') '6.62'
Set-CellText $ws.Range('E26') '  -3.57%  '
Set-CellText $ws.Range('D27') '94.62'
Set-CellText $ws.Range('E27') '  -1.12%  '
Set-CellText $ws.Range('B28') 'Aptos'
Set-CellText $ws.Range('C28') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-CellText $ws.Range('D28') '12.02'
Set-CellText $ws.Range('E28') '  -0.83%  '
Set-CellText $ws.Range('B29') 'WrappedeETH'
Set-CellText $ws.Range('C29') 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-CellText $ws.Range('D29') '3.627.32'
Set-CellText $ws.Range('E29') '  +2.99%  '
Set-CellText $ws.Range('D30') '11.85'
Set-CellText $ws.Range('E30') '  +5.32%  '
Set-CellText $ws.Range('B31') 'PancakeSwap'
Set-CellText $ws.Range('C31') 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-CellText $ws.Range('D31') '2.80'
Set-CellText $ws.Range('E31') '  +9.98%  '
Set-CellText $ws.Range('B32') 'Dai'
Set-CellText $ws.Range('C32') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellText $ws.Range('D32') '1.00'
Set-CellText $ws.Range('E32') '  +0.13%  '
Set-CellText $ws.Range('D33') '0.139'
Set-CellText $ws.Range('E33') '  -1.66%  '
Set-CellText $ws.Range('B34') 'Binance-PegBSC-USD'
Set-CellText $ws.Range('C34') 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-CellText $ws.Range('D34') '1.00'
Set-CellText $ws.Range('E34') '  +0.08%  '
Set-CellText $ws.Range('B35') 'Cronos'
Set-CellText $ws.Range('C35') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText $ws.Range('D35') '0.178'
Set-CellText $ws.Range('E35') '  -3.76%  '
Set-CellText $ws.Range('D36') '30.16'
Set-CellText $ws.Range('E36') '  +6.69%  '
Set-CellText $ws.Range('D37') '0.558'
Set-CellText $ws.Range('E37') '  +0.94%  '
Set-CellText $ws.Range('D38') '557.40'
Set-CellText $ws.Range('E38') '  +5.16%  '
Set-CellText $ws.Range('D39') '7.69'
Set-CellText $ws.Range('E39') '  -3.79%  '
Set-CellText $ws.Range('D40') '1.46'
Set-CellText $ws.Range('E40') '  -3.12%  '
Set-CellText $ws.Range('E41') '  +1.04%  '
Set-CellText $ws.Range('E42') '  -0.03%  '
Set-CellText $ws.Range('D43') '0.912'
Set-CellText $ws.Range('E43') '  +8.84%  '
Set-CellText $ws.Range('D44') '24.05'
Set-CellText $ws.Range('E44') '  -1.37%  '
Set-CellText $ws.Range('D45') '1.73'
Set-CellText $ws.Range('E45') '  +0.70%  '
Set-CellText $ws.Range('B46') 'Filecoin'
Set-CellText $ws.Range('C46') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws.Range('D46') '5.63'
Set-CellText $ws.Range('E46') '  +2.40%  '
Set-CellText $ws.Range('B47') 'MantraDAO'
Set-CellText $ws.Range('C47') 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-CellText $ws.Range('D47') '3.61'
Set-CellText $ws.Range('E47') '  +0.28%  '
Set-CellText $ws.Range('B48') 'VeChain'
Set-CellText $ws.Range('C48') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws.Range('D48') '0.0412'
Set-CellText $ws.Range('E48') '  -2.59%  '
Set-CellText $ws.Range('B49') 'dogwifhat'
Set-CellText $ws.Range('C49') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText $ws.Range('D49') '3.33'
Set-CellText $ws.Range('E49') '  +5.45%  '
Set-CellText $ws.Range('B50') 'Stacks'
Set-CellText $ws.Range('C50') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CellText $ws.Range('D50') '2.19'
Set-CellText $ws.Range('E50') '  -3.79%  '
Set-CellText $ws.Range('D51') '53.82'
Set-CellText $ws.Range('E51') '  -1.47%  '
